$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.008627148837446561
$ws.Range("J2").Value = 0.008627148837446563
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.119639
$ws.Range("N2").Value = 0.358917
$ws.Range("O2").Value = 0.02933906125009379
$ws.Range("P2").Value = 0.03107835843382883
$ws.Range("Q2").Value = 0.002077331836666666
$ws.Range("R2").Value = 0.01869598653
$ws.Range("S2").Value = 0.0002531124481555201
$ws.Range("T2").Value = 0.000268117623832154
$ws.Range("I3").Value = 0.008627148837446561
$ws.Range("J3").Value = 0.008627148837446563
$ws.Range("O3").Value = 0.05256748359289284
$ws.Range("P3").Value = 0.05568382311683939
$ws.Range("S3").Value = 0.0004535075049659166
$ws.Range("T3").Value = 0.000480392629867021
$ws.Range("I4").Value = 0.008627148837446561
$ws.Range("J4").Value = 0.008627148837446563
$ws.Range("M4").Value = 1.018537666666667
$ws.Range("N4").Value = 3.055613
$ws.Range("O4").Value = 0.2497759007335481
$ws.Range("P4").Value = 0.2645832770503125
$ws.Range("Q4").Value = 0.01768520901888889
$ws.Range("R4").Value = 0.15916688117
$ws.Range("S4").Value = 0.002154853871635597
$ws.Range("T4").Value = 0.002282599311012406
$ws.Range("I5").Value = 0.008627148837446561
$ws.Range("J5").Value = 0.008627148837446563
$ws.Range("M5").Value = 0.684642
$ws.Range("N5").Value = 1.369284
$ws.Range("O5").Value = 0.1678946963146358
$ws.Range("P5").Value = 0.1185652921140734
$ws.Range("Q5").Value = 0.01188766726
$ws.Range("R5").Value = 0.07132600355999999
$ws.Range("S5").Value = 0.001448452534124254
$ws.Range("T5").Value = 0.00102288042202344
$ws.Range("I6").Value = 0.008627148837446561
$ws.Range("J6").Value = 0.008627148837446563
$ws.Range("M6").Value = 2.040627333333334
$ws.Range("N6").Value = 6.121882
$ws.Range("O6").Value = 0.5004228581088294
$ws.Range("P6").Value = 0.5300892492849458
$ws.Range("Q6").Value = 0.03543209259777778
$ws.Range("R6").Value = 0.31888883338
$ws.Range("S6").Value = 0.004317222478565272
$ws.Range("T6").Value = 0.004573158850711542
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.023931333333333
$ws.Range("H7").Value = 3.071794
$ws.Range("I7").Value = 0.5087507014009469
$ws.Range("J7").Value = 0.5087507014009469
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.119639
$ws.Range("N7").Value = 0.358917
$ws.Range("O7").Value = 0.02933906125009379
$ws.Range("P7").Value = 0.03107835843382883
$ws.Range("Q7").Value = 0.1225021207886667
$ws.Range("R7").Value = 1.102519087098
$ws.Range("S7").Value = 0.01492626798943056
$ws.Range("T7").Value = 0.01581113665160045
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.023931333333333
$ws.Range("H8").Value = 3.071794
$ws.Range("I8").Value = 0.5087507014009469
$ws.Range("J8").Value = 0.5087507014009469
$ws.Range("O8").Value = 0.05256748359289284
$ws.Range("P8").Value = 0.05568382311683939
$ws.Range("Q8").Value = 0.2194899206133333
$ws.Range("R8").Value = 1.97540928552
$ws.Range("S8").Value = 0.026743744148767
$ws.Range("T8").Value = 0.0283291840673783
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.023931333333333
$ws.Range("H9").Value = 3.071794
$ws.Range("I9").Value = 0.5087507014009469
$ws.Range("J9").Value = 0.5087507014009469
$ws.Range("M9").Value = 1.018537666666667
$ws.Range("N9").Value = 3.055613
$ws.Range("O9").Value = 0.2497759007335481
$ws.Range("P9").Value = 0.2645832770503125
$ws.Range("Q9").Value = 1.042912631080222
$ws.Range("R9").Value = 9.386213679721999
$ws.Range("S9").Value = 0.1270736646912459
$ws.Range("T9").Value = 0.1346069277783076
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.023931333333333
$ws.Range("H10").Value = 3.071794
$ws.Range("I10").Value = 0.5087507014009469
$ws.Range("J10").Value = 0.5087507014009469
$ws.Range("M10").Value = 0.684642
$ws.Range("N10").Value = 1.369284
$ws.Range("O10").Value = 0.1678946963146358
$ws.Range("P10").Value = 0.1185652921140734
$ws.Range("Q10").Value = 0.7010263959159999
$ws.Range("R10").Value = 4.206158375495999
$ws.Range("S10").Value = 0.08541654451156995
$ws.Range("T10").Value = 0.06032017552484299
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.023931333333333
$ws.Range("H11").Value = 3.071794
$ws.Range("I11").Value = 0.5087507014009469
$ws.Range("J11").Value = 0.5087507014009469
$ws.Range("M11").Value = 2.040627333333334
$ws.Range("N11").Value = 6.121882
$ws.Range("O11").Value = 0.5004228581088294
$ws.Range("P11").Value = 0.5300892492849458
$ws.Range("Q11").Value = 2.089462266256445
$ws.Range("R11").Value = 18.805160396308
$ws.Range("S11").Value = 0.2545904800599335
$ws.Range("T11").Value = 0.2696832773788176
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.9713440000000001
$ws.Range("H12").Value = 2.914032
$ws.Range("I12").Value = 0.4826221497616066
$ws.Range("J12").Value = 0.4826221497616065
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.119639
$ws.Range("N12").Value = 0.358917
$ws.Range("O12").Value = 0.02933906125009379
$ws.Range("P12").Value = 0.03107835843382883
$ws.Range("Q12").Value = 0.116210624816
$ws.Range("R12").Value = 1.045895623344
$ws.Range("S12").Value = 0.01415968081250771
$ws.Range("T12").Value = 0.01499910415839622
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.9713440000000001
$ws.Range("H13").Value = 2.914032
$ws.Range("I13").Value = 0.4826221497616066
$ws.Range("J13").Value = 0.4826221497616065
$ws.Range("O13").Value = 0.05256748359289284
$ws.Range("P13").Value = 0.05568382311683939
$ws.Range("Q13").Value = 0.20821729984
$ws.Range("R13").Value = 1.87395569856
$ws.Range("S13").Value = 0.02537023193915992
$ws.Range("T13").Value = 0.02687424641959407
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.9713440000000001
$ws.Range("H14").Value = 2.914032
$ws.Range("I14").Value = 0.4826221497616066
$ws.Range("J14").Value = 0.4826221497616065
$ws.Range("M14").Value = 1.018537666666667
$ws.Range("N14").Value = 3.055613
$ws.Range("O14").Value = 0.2497759007335481
$ws.Range("P14").Value = 0.2645832770503125
$ws.Range("Q14").Value = 0.9893504512906669
$ws.Range("R14").Value = 8.904154061616001
$ws.Range("S14").Value = 0.1205473821706666
$ws.Range("T14").Value = 0.1276937499609926
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.9713440000000001
$ws.Range("H15").Value = 2.914032
$ws.Range("I15").Value = 0.4826221497616066
$ws.Range("J15").Value = 0.4826221497616065
$ws.Range("M15").Value = 0.684642
$ws.Range("N15").Value = 1.369284
$ws.Range("O15").Value = 0.1678946963146358
$ws.Range("P15").Value = 0.1185652921140734
$ws.Range("Q15").Value = 0.665022898848
$ws.Range("R15").Value = 3.990137393088
$ws.Range("S15").Value = 0.08102969926894163
$ws.Range("T15").Value = 0.05722223616720695
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.9713440000000001
$ws.Range("H16").Value = 2.914032
$ws.Range("I16").Value = 0.4826221497616066
$ws.Range("J16").Value = 0.4826221497616065
$ws.Range("M16").Value = 2.040627333333334
$ws.Range("N16").Value = 6.121882
$ws.Range("O16").Value = 0.5004228581088294
$ws.Range("P16").Value = 0.5300892492849458
$ws.Range("Q16").Value = 1.982151116469334
$ws.Range("R16").Value = 17.839360048224
$ws.Range("S16").Value = 0.2415151555703306
$ws.Range("T16").Value = 0.2558328130554167
